$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New job posting row (row 13): JD_012 / Junior Developer
# Job_Description reuses the same text as the "Software Engineer" posting (row 2, column C)
$jobDescription = $ws.Range("C2").Value2

$ws.Range("A13").Value = "JD_012"
$ws.Range("B13").Value = "Junior Developer"
$ws.Range("C13").Value = $jobDescription
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 5

# Setting a long, wrapped description can cause Excel to mark the row with an
# explicit custom height; AutoFit restores the default (non-custom) row height.
$ws.Rows.Item(13).AutoFit() | Out-Null
